$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = '{''model__num_layers'': 2, ''model__num_epochs'': 100, ''model__learning_rate'': 0.001, ''model__hidden_size'': 100}'
$ws.Range("D2").Value = -172.8527694896686
$ws.Range("E2").Value = -261.2691374453992
$ws.Range("F2").Value = -0.1555368039808629
$ws.Range("G2").Value = 15.72131944615002
$ws.Range("H2").Value = 691.3090721403909
$ws.Range("R2").Value = 144.8301923274994

# Row 3
$ws.Range("C3").Value = '{''model__num_layers'': 2, ''model__num_epochs'': 1000, ''model__learning_rate'': 0.01, ''model__hidden_size'': 50}'
$ws.Range("D3").Value = -76.6727051768876
$ws.Range("E3").Value = -242.2800866981888
$ws.Range("F3").Value = -1.37409088851563
$ws.Range("G3").Value = 24.9399623209483
$ws.Range("H3").Value = 1170.379548839687
$ws.Range("R3").Value = 158.2019264698029

# Row 4
$ws.Range("C4").Value = '{''model__num_layers'': 1, ''model__num_epochs'': 100, ''model__learning_rate'': 0.001, ''model__hidden_size'': 100}'
$ws.Range("D4").Value = -196.1482591672423
$ws.Range("E4").Value = -298.8487153428299
$ws.Range("F4").Value = 0.002907464150800033
$ws.Range("G4").Value = 15.05971658410547
$ws.Range("H4").Value = 628.0208864998216
$ws.Range("R4").Value = 151.3259847164154

# Row 5
$ws.Range("D5").Value = -161.7739497669316
$ws.Range("E5").Value = -227.0103812347398
$ws.Range("F5").Value = -0.01199033354670448
$ws.Range("G5").Value = 14.97530084975778
$ws.Range("H5").Value = 627.1710922436423
$ws.Range("R5").Value = 81.75265431404114

# Row 6
$ws.Range("C6").Value = '{''model__num_layers'': 2, ''model__num_epochs'': 100, ''model__learning_rate'': 0.001, ''model__hidden_size'': 100}'
$ws.Range("D6").Value = -152.0897618055957
$ws.Range("E6").Value = -252.3588980430893
$ws.Range("F6").Value = -0.07467443601760616
$ws.Range("G6").Value = 15.20127800723733
$ws.Range("H6").Value = 659.409719107261
$ws.Range("R6").Value = 74.70702195167542

# Row 7
$ws.Range("C7").Value = '{''model__num_layers'': 1, ''model__num_epochs'': 200, ''model__learning_rate'': 0.01, ''model__hidden_size'': 100}'
$ws.Range("D7").Value = -135.7946958025966
$ws.Range("E7").Value = -212.0371116356648
$ws.Range("F7").Value = -0.1371438836669445
$ws.Range("G7").Value = 16.86700426136522
$ws.Range("H7").Value = 649.4439889931624
$ws.Range("R7").Value = 67.73950147628784

# Row 8
$ws.Range("C8").Value = '{''model__num_layers'': 1, ''model__num_epochs'': 200, ''model__learning_rate'': 0.001, ''model__hidden_size'': 50}'
$ws.Range("D8").Value = -135.7096714671342
$ws.Range("E8").Value = -282.0976000018738
$ws.Range("F8").Value = -0.0616693587328205
$ws.Range("G8").Value = 15.32435088053141
$ws.Range("H8").Value = 641.645778393803
$ws.Range("R8").Value = 127.9376130104065

# Row 9
$ws.Range("C9").Value = '{''model__num_layers'': 1, ''model__num_epochs'': 100, ''model__learning_rate'': 0.01, ''model__hidden_size'': 50}'
$ws.Range("D9").Value = -170.3885490090783
$ws.Range("E9").Value = -229.3575651894215
$ws.Range("F9").Value = -0.02148651825069288
$ws.Range("G9").Value = 14.98958668103501
$ws.Range("H9").Value = 647.2581252329063
$ws.Range("R9").Value = 76.35982370376587

# Row 10
$ws.Range("D10").Value = -130.1377172414571
$ws.Range("E10").Value = -259.9469519184787
$ws.Range("F10").Value = -0.02723020329831707
$ws.Range("G10").Value = 15.09695897312469
$ws.Range("H10").Value = 625.1326716294051
$ws.Range("R10").Value = 8.293279647827148

# Row 11
$ws.Range("D11").Value = -171.4621561722359
$ws.Range("E11").Value = -245.6482793673779
$ws.Range("F11").Value = -0.03677607512387978
$ws.Range("G11").Value = 15.50050531474423
$ws.Range("H11").Value = 636.0944931596936
$ws.Range("R11").Value = 38.15863180160522

# Row 12
$ws.Range("C12").Value = '{''model__num_layers'': 2, ''model__num_epochs'': 100, ''model__learning_rate'': 0.001, ''model__hidden_size'': 100}'
$ws.Range("D12").Value = -165.6760429019877
$ws.Range("E12").Value = -237.7659667535859
$ws.Range("F12").Value = -0.02446209747770381
$ws.Range("G12").Value = 15.11521404488446
$ws.Range("H12").Value = 657.4556169624924
$ws.Range("R12").Value = 57.44451189041138

# Row 13
$ws.Range("D13").Value = -178.9240506885007
$ws.Range("E13").Value = -221.0614181646161
$ws.Range("F13").Value = -0.03471829223772793
$ws.Range("G13").Value = 15.10966319262709
$ws.Range("H13").Value = 655.9333324449693
$ws.Range("R13").Value = 50.47422766685486
